# Auto-generated script to apply 2022-09-12 daily crime-data update
# to output/violent-crime-full-year.xlsx (column I = year 2022 running totals).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 5062
$ws.Range('I3').Value = 5274
$ws.Range('I4').Value = 1208
$ws.Range('I6').Value = 5778
$ws.Range('I7').Value = 17811

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 179
$ws.Range('I3').Value = 185
$ws.Range('I6').Value = 160
$ws.Range('I7').Value = 572

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I3').Value = 119
$ws.Range('I4').Value = 23
$ws.Range('I6').Value = 86
$ws.Range('I7').Value = 331

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I3').Value = 251
$ws.Range('I6').Value = 224
$ws.Range('I7').Value = 697

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I3').Value = 47
$ws.Range('I7').Value = 147

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I7').Value = 576
$ws.Range('I8').Value = 1075
$ws.Range('I13').Value = 32
$ws.Range('I18').Value = 130
$ws.Range('I19').Value = 488
$ws.Range('I23').Value = 173
$ws.Range('I25').Value = 89
$ws.Range('I26').Value = 25
$ws.Range('I27').Value = 166
$ws.Range('I29').Value = 1131
$ws.Range('I32').Value = 24
$ws.Range('I33').Value = 813
$ws.Range('I36').Value = 233
$ws.Range('I37').Value = 572
$ws.Range('I42').Value = 604
$ws.Range('I43').Value = 140
$ws.Range('I50').Value = 81
$ws.Range('I51').Value = 192
$ws.Range('I52').Value = 395
$ws.Range('I54').Value = 383
$ws.Range('I63').Value = 62
$ws.Range('I67').Value = 697
$ws.Range('I68').Value = 66
$ws.Range('I78').Value = 253
$ws.Range('I79').Value = 511
$ws.Range('I83').Value = 375
$ws.Range('I84').Value = 147
$ws.Range('I85').Value = 803
$ws.Range('I87').Value = 36
$ws.Range('I88').Value = 163
$ws.Range('I93').Value = 103
$ws.Range('I94').Value = 178
$ws.Range('I97').Value = 137
$ws.Range('I98').Value = 123
$ws.Range('I99').Value = 331
$ws.Range('I101').Value = 17811

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I3').Value = 141
$ws.Range('I6').Value = 71
$ws.Range('I7').Value = 375

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I2').Value = 185
$ws.Range('I7').Value = 813

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I3').Value = 85
$ws.Range('I7').Value = 383

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 337
$ws.Range('I3').Value = 391
$ws.Range('I6').Value = 305
$ws.Range('I7').Value = 1131

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I3').Value = 147
$ws.Range('I7').Value = 488

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I6').Value = 195
$ws.Range('I7').Value = 803

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I3').Value = 207
$ws.Range('I7').Value = 604

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('I4').Value = 10
$ws.Range('I6').Value = 32

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I3').Value = 63
$ws.Range('I7').Value = 253

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I2').Value = 49
$ws.Range('I7').Value = 173

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I2').Value = 146
$ws.Range('I3').Value = 165
$ws.Range('I6').Value = 151
$ws.Range('I7').Value = 511

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('I6').Value = 55
$ws.Range('I7').Value = 130

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I6').Value = 73
$ws.Range('I7').Value = 233

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('I6').Value = 43
$ws.Range('I7').Value = 103

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I2').Value = 107
$ws.Range('I6').Value = 98
$ws.Range('I7').Value = 395

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 103
$ws.Range('I7').Value = 178

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('I3').Value = 27
$ws.Range('I7').Value = 89

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('I6').Value = 82
$ws.Range('I7').Value = 123

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('I4').Value = 18
$ws.Range('I7').Value = 81

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('I6').Value = 15
$ws.Range('I7').Value = 25

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I6').Value = 83
$ws.Range('I7').Value = 137

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I3').Value = 58
$ws.Range('I6').Value = 46
$ws.Range('I7').Value = 163

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('I2').Value = 9
$ws.Range('I7').Value = 24

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 337
$ws.Range('I3').Value = 303
$ws.Range('I6').Value = 345
$ws.Range('I7').Value = 1075

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('I2').Value = 47
$ws.Range('I3').Value = 31
$ws.Range('I7').Value = 166

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I6').Value = 76
$ws.Range('I7').Value = 192

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('I3').Value = 21
$ws.Range('I7').Value = 66

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I2').Value = 28
$ws.Range('I7').Value = 140

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I3').Value = 178
$ws.Range('I7').Value = 576

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('I2').Value = 4
$ws.Range('I7').Value = 36
